$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Veranstaltungsnummer" (column D) and "Slot" (column E) values
# for randomly-assigned ("Zufaellig") rows, per the source diff.
$updates = @(
    @{ Row = 6; D = 12; E = $null }
    @{ Row = 11; D = 8; E = $null }
    @{ Row = 16; D = 13; E = $null }
    @{ Row = 21; D = 7; E = $null }
    @{ Row = 26; D = 13; E = $null }
    @{ Row = 36; D = 14; E = $null }
    @{ Row = 41; D = 27; E = $null }
    @{ Row = 56; D = 16; E = $null }
    @{ Row = 80; D = 19; E = 3 }
    @{ Row = 81; D = 12; E = 5 }
    @{ Row = 86; D = 7; E = $null }
    @{ Row = 91; D = 1; E = $null }
    @{ Row = 96; D = 14; E = $null }
    @{ Row = 101; D = 1; E = $null }
    @{ Row = 116; D = 3; E = $null }
    @{ Row = 121; D = 7; E = $null }
    @{ Row = 126; D = 13; E = $null }
    @{ Row = 131; D = 11; E = $null }
    @{ Row = 141; D = 1; E = $null }
    @{ Row = 146; D = 14; E = $null }
    @{ Row = 150; D = 5; E = $null }
    @{ Row = 151; D = 4; E = $null }
    @{ Row = 171; D = 5; E = $null }
    @{ Row = 176; D = 2; E = $null }
    @{ Row = 181; D = 3; E = $null }
    @{ Row = 200; D = 15; E = 2 }
    @{ Row = 201; D = 27; E = 5 }
    @{ Row = 206; D = 3; E = $null }
    @{ Row = 209; D = 5; E = $null }
    @{ Row = 210; D = 16; E = 5 }
    @{ Row = 211; D = 15; E = 2 }
    @{ Row = 231; D = 7; E = $null }
    @{ Row = 241; D = 13; E = $null }
    @{ Row = 256; D = 15; E = $null }
    @{ Row = 265; D = 1; E = $null }
    @{ Row = 266; D = 10; E = $null }
    @{ Row = 271; D = 27; E = $null }
    @{ Row = 276; D = 10; E = $null }
    @{ Row = 281; D = 12; E = $null }
    @{ Row = 285; D = 14; E = 3 }
    @{ Row = 286; D = 1; E = 4 }
    @{ Row = 291; D = 1; E = $null }
    @{ Row = 296; D = 12; E = $null }
    @{ Row = 306; D = 5; E = $null }
    @{ Row = 356; D = 21; E = $null }
    @{ Row = 366; D = 12; E = $null }
    @{ Row = 375; D = 21; E = 2 }
    @{ Row = 376; D = 3; E = 4 }
    @{ Row = 391; D = 16; E = $null }
    @{ Row = 396; D = 26; E = $null }
    @{ Row = 406; D = 10; E = $null }
    @{ Row = 421; D = 10; E = $null }
    @{ Row = 426; D = 4; E = $null }
    @{ Row = 431; D = 19; E = $null }
    @{ Row = 436; D = 7; E = $null }
    @{ Row = 441; D = 12; E = $null }
    @{ Row = 466; D = 9; E = $null }
    @{ Row = 476; D = 26; E = $null }
    @{ Row = 486; D = 10; E = $null }
    @{ Row = 491; D = 2; E = $null }
    @{ Row = 496; D = 16; E = $null }
    @{ Row = 505; D = 17; E = 4 }
    @{ Row = 506; D = 10; E = 2 }
    @{ Row = 510; D = 6; E = $null }
    @{ Row = 511; D = 9; E = $null }
    @{ Row = 516; D = 12; E = $null }
    @{ Row = 520; D = 26; E = $null }
    @{ Row = 521; D = 6; E = $null }
    @{ Row = 546; D = 4; E = $null }
    @{ Row = 556; D = 3; E = $null }
    @{ Row = 561; D = 16; E = $null }
    @{ Row = 575; D = 1; E = 4 }
    @{ Row = 576; D = 27; E = 5 }
    @{ Row = 581; D = 1; E = $null }
    @{ Row = 586; D = 6; E = $null }
    @{ Row = 591; D = 22; E = $null }
    @{ Row = 616; D = 16; E = $null }
    @{ Row = 621; D = 16; E = $null }
    @{ Row = 626; D = 1; E = $null }
    @{ Row = 631; D = 6; E = $null }
    @{ Row = 636; D = 13; E = $null }
    @{ Row = 641; D = 15; E = $null }
    @{ Row = 646; D = 13; E = $null }
    @{ Row = 650; D = $null; E = 4 }
    @{ Row = 651; D = 6; E = 5 }
    @{ Row = 656; D = 22; E = $null }
    @{ Row = 661; D = 15; E = $null }
    @{ Row = 666; D = 20; E = $null }
    @{ Row = 667; D = 17; E = 4 }
    @{ Row = 668; D = 10; E = 2 }
    @{ Row = 669; D = $null; E = 3 }
    @{ Row = 671; D = 23; E = 1 }
    @{ Row = 672; D = 26; E = 1 }
    @{ Row = 673; D = 13; E = 5 }
    @{ Row = 674; D = 25; E = 2 }
    @{ Row = 675; D = 1; E = 4 }
    @{ Row = 676; D = 15; E = 3 }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) { $ws.Cells.Item($u.Row, 4).Value = $u.D }
    if ($null -ne $u.E) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
}
